$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the note text in D2 (shared string index 6)
$ws.Range("D2").Value = "Sv_0699 was halved with extract left in vial."

# 2. Fill in C5 with "T" (same as C2)
$ws.Range("C5").Value = "T"

# 3. Fill in D5 with the rich-text note about Sv_0158/Sv_0104, two runs with
#    different formatting (second run uses a plain Arial/10 font, no colour)
$ws.Range("D5").Value = "Sv_0158 and Sv_0104 were halved with extract left in vial."
$ws.Range("D5").Characters(26, 33).Font.ColorIndex = -4105
$ws.Range("D5").Characters(26, 33).Font.Name = "Arial"
$ws.Range("D5").Characters(26, 33).Font.Size = 10

# 4. Wrap text + taller row for row 5, to show the new note
$ws.Range("D5").WrapText = $true
$ws.Rows(5).RowHeight = 24.85

# 5. Add a new merged row-group (rows 8-10) mirroring rows 5-7 / 2-4
#    Copy the whole A:C block in one go so existing styles (s=4 / s=5) get
#    reused instead of new (duplicate) cellXfs being generated.
$ws.Range("A6:C7").Copy()
$ws.Range("A9:C10").PasteSpecial(-4122)
$ws.Range("A6:C6").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

$ws.Range("A8").Value = "0001-0100"
$ws.Range("C8").Value = "T"

$ws.Range("A8:A10").Merge()
$ws.Range("B8:B10").Merge()
$ws.Range("C8:C10").Merge()

# 6. Move selection to A18, like in the final workbook
$ws.Range("A18").Select()
